$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new document-metadata columns for row 2 (category, title, accom weight, study weight)
$ws.Range("D2").Value = "student"
$ws.Range("E2").Value = "Cost of Living in Poland"
$ws.Range("H2").Value = 2
$ws.Range("J2").Value = 20

# Move the active selection to J2, matching the updated view state
$ws.Range("J2").Select()
